# Applies the commit's changes:
#  - Metadata!B8: updated Date value
#  - Elements sheet: several Binding Value Set URLs (CodeSystem instead of ValueSet)
#  - Elements sheet: Type(s) column for telecommunication/adresseSE cells (drop
#    the "ContactPoint {...}" / "Address {...}" wrapper, keep just the URL)
#  - Column widths for columns K and Z adjusted (best-fit) to match new content

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- Metadata: Date ---
$wsMetadata.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Elements: Binding Value Set URLs ---
$wsElements.Range("Z5").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R23-ModeExercice?vs"
$wsElements.Range("Z9").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R22-GenreActivite?vs"
$wsElements.Range("Z10").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R25-MotifFinActivite?vs"
$wsElements.Range("Z11").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R32-StatutHospitalier?vs"
$wsElements.Range("Z12").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R06-SectionTableauCNOP?vs"
$wsElements.Range("Z13").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-G05-SousSectionTableauCNOP?vs"
$wsElements.Range("Z14").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R24-TypeActiviteLiberale?vs"
$wsElements.Range("Z15").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R34-StatutProfessionnelSSA?vs"

# --- Elements: Type(s) column, drop the wrapper around the URL ---
$wsElements.Range("K16").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/Telecommunication`n"
$wsElements.Range("K17").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/Adresse`n"

# --- Column widths: re-fit columns K (11) and Z (26) to their new best content width ---
# (values chosen so the engine's pixel-grid rounding lands on the closest
# representable width to the authored target: K -> 64.53125, Z -> 76.68359375)
$wsElements.Columns.Item(11).ColumnWidth = 63.65
$wsElements.Columns.Item(26).ColumnWidth = 75.75
